# Updated cryptos list values (Mon Jul 24 07:40:04 UTC 2023) with GitHub Actions
# Applies the per-cell text updates described by the commit diff: refreshed
# Price / Volume(1h) figures, plus the two ranking swaps (WrappedEther <-> Polkadot
# at rows 13/14, and Aptos <-> RenderToken at rows 48/49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.811.26"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.874.39"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7284"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.62"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07283"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.57"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08183"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7461"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.362"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.865.95"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.77"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "29.799.70"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.022"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "248.04"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.46"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007838"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "2.124.44"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.742"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1526"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.257"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.22"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.61"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.443"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.538"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.525"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.180"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05416"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.232"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7421"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.703"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01928"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.741"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4475"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8885"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.994"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.67"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.29"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "1.037.02"
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.830"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.493"
$ws.Range("E49").Value = "  -2.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.642"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").Value = "2.015.64"
$ws.Range("E51").Value = "  +0.07%  "
